# Update column A on Sheet1 so the per-month trading-day index (which
# restarted at 0 every month) becomes one continuous running index for
# the whole year: A{r} = r - 2 for every data row r (rows 2..245).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 245
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
